$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.366.55'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '3.933.51'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''487.16'
$ws.Range('E5').Value = '  +3.50%  '
$ws.Range('D6').Value = '''148.43'
$ws.Range('E6').Value = '  +2.34%  '
$ws.Range('D7').Value = '''0.628'
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  +1.30%  '
$ws.Range('E10').Value = '  +4.27%  '
$ws.Range('E11').Value = '  +5.49%  '
$ws.Range('D12').Value = '''43.11'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').Value = '''10.72'
$ws.Range('E13').Value = '  +3.73%  '
$ws.Range('D14').Value = '4.560.13'
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').Value = '''14.62'
$ws.Range('E15').Value = '  -2.13%  '
$ws.Range('D16').Value = '3.918.16'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').Value = '''0.137'
$ws.Range('E17').Value = '  -0.51%  '
$ws.Range('D18').Value = '''20.04'
$ws.Range('E18').Value = '  +1.27%  '
$ws.Range('E19').Value = '  -1.07%  '
$ws.Range('D20').Value = '68.451.03'
$ws.Range('D21').Value = '''444.06'
$ws.Range('E21').Value = '  +3.29%  '
$ws.Range('D22').Value = '''3.51'
$ws.Range('E22').Value = '  +4.46%  '
$ws.Range('D23').Value = '''15.21'
$ws.Range('E23').Value = '  +4.84%  '
$ws.Range('D24').Value = '''88.62'
$ws.Range('E24').Value = '  +1.37%  '
$ws.Range('D25').Value = '''11.49'
$ws.Range('E25').Value = '  +20.13%  '
$ws.Range('D26').Value = '''11.50'
$ws.Range('E26').Value = '  +13.03%  '
$ws.Range('D27').Value = '''3.65'
$ws.Range('E27').Value = '  +2.80%  '
$ws.Range('D28').Value = '''39.03'
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('E29').Value = '  +1.92%  '
$ws.Range('D30').Value = '''724.00'
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').Value = '''13.75'
$ws.Range('E31').Value = '  +1.75%  '
$ws.Range('E32').Value = '  -0.97%  '
$ws.Range('D33').Value = '''2.92'
$ws.Range('E33').Value = '  +4.22%  '
$ws.Range('D34').Value = '0.0₃0916'
$ws.Range('E34').Value = '  +16.82%  '
$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = '''42.63'
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = '''6.26'
$ws.Range('E36').Value = '  +16.78%  '
$ws.Range('D37').Value = '''61.22'
$ws.Range('E37').Value = '  +5.79%  '
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('D39').Value = '''0.402'
$ws.Range('E39').Value = '  +19.71%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '''3.05'
$ws.Range('E40').Value = '  +18.25%  '
$ws.Range('B41').Value = 'Dai'
$ws.Range('C41').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('E42').Value = '  +7.51%  '
$ws.Range('D43').Value = '''0.0483'
$ws.Range('E43').Value = '  +1.53%  '
$ws.Range('E44').Value = '  +5.36%  '
$ws.Range('E45').Value = '  +1.73%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '''1.00'
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('B47').Value = 'LidoDAOToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D47').Value = '''3.43'
$ws.Range('E47').Value = '  +0.55%  '
$ws.Range('B48').Value = 'ApeXProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D48').Value = '''3.27'
$ws.Range('E48').Value = '  +3.34%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0₆0355'
$ws.Range('E49').Value = '  +39.69%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').Value = '''2.16'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').Value = '''145.67'
$ws.Range('E51').Value = '  -0.44%  '
